# Commit: "add support for hitPolicy and annotations"
#
# The example DMN decision table gains an annotation column: row 3
# (amount "> 800" / invoiceCategory "Travel Expenses" / result "sales")
# gets a new "anno1" value in column D.
#
# We copy the formatting of the adjoining C3 cell (same style index as the
# rest of the table) with PasteSpecial(xlPasteFormats = -4122) before
# writing the value, so D3 ends up styled like its row-mates instead of
# picking up a brand-new, unstyled cell format.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("C3").Copy()
$ws.Range("D3").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("D3").Value = "anno1"
